$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column values are numeric-looking text (e.g. "1.003", "27.586.29") that must
# stay as text, matching the source inlineStr cells. Force text format before
# assigning so Excel does not auto-convert to a number, then reset the style back
# to Normal so no stray style index is left attached to the cell.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.586.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.53%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.750.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.29%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.28%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4589"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +9.69%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3576"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.52%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07479"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.03%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.79%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.091"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.72%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.01%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.97%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.998"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.52%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.081"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.33%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.748.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.02%  "

# Row 17
$ws.Range("E17").Value = "  +1.93%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001063"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.55%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06433"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.77%  "

# Row 20
$ws.Range("E20").Value = "  -0.02%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.16%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.801"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.94%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.648.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.43%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.78%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.109"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.02%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.95%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.74%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.952.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.87%  "

# Row 29
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.063"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.95%  "

# Row 30
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.20%  "

# Row 31
$ws.Range("E31").Value = "  -6.03%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09213"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.36%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.671"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.82%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.520"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.20%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.80"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.63%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02293"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.24%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06040"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.09%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2089"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.13%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.968"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.74%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6304"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.59%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.204"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.12%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.378"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.02%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.764"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.03%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.80%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5897"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.99%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.716"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.43%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.54%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.937"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.07%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.132"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.12%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06860"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.98%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.89%  "
